# Add a "Turkey" test-data sheet (copy of the "Spain" sheet layout) with
# Turkey-specific content, mirroring the other per-country sheets in this
# workbook ("Added Test data for Turkey Template for Zettler").

$wb = $excel.ActiveWorkbook

$spain = $wb.Worksheets.Item("Spain")

# Duplicate the Spain sheet (this carries over styles, merged cells,
# column widths/bestFit flags, page setup, etc.) and place it right after
# Spain - i.e. as the new last sheet.
$spain.Copy($null, $spain)
$turkey = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Name = "Turkey"

# Update the market name and ticket reference for Turkey.
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3291/T3300 "

# The Turkey text is shorter than Spain's, so it no longer wraps to two
# lines: rows 3 and 5 go back to the default height, row 4 (bigger font)
# settles at 15.6 instead of Spain's wrapped 28.8.
$turkey.Rows(3).RowHeight = 14.4
$turkey.Rows(4).RowHeight = 15.6
$turkey.Rows(5).RowHeight = 14.4

# Columns resize to fit the new (shorter) Turkey content.
$turkey.Columns("B").ColumnWidth = 23.166666666666668
$turkey.Columns("D").ColumnWidth = 19.666666666666668

# Selections: Spain is no longer the active tab, and its selection
# reverts to the full used range; Turkey becomes the active tab with its
# selection parked at F12.
$spain.Range("A1:D11").Select() | Out-Null
$turkey.Range("F12").Select() | Out-Null

$turkey.Activate()
